$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.9584923333333334
$ws.Range("H2").Value = 2.875477
$ws.Range("I2").Value = 0.2532195598902293
$ws.Range("J2").Value = 0.2532195598902293
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 5.399171333333334
$ws.Range("N2").Value = 16.197514
$ws.Range("O2").Value = 0.1822675802569684
$ws.Range("P2").Value = 0.1822675802569683
$ws.Range("Q2").Value = 5.175064329353113
$ws.Range("R2").Value = 46.57557896417801
$ws.Range("S2").Value = 0.04615371645492657
$ws.Range("T2").Value = 0.04615371645492657

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.9584923333333334
$ws.Range("H3").Value = 2.875477
$ws.Range("I3").Value = 0.2532195598902293
$ws.Range("J3").Value = 0.2532195598902293
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.704475666666667
$ws.Range("N3").Value = 11.113427
$ws.Range("O3").Value = 0.1250573049452193
$ws.Range("P3").Value = 0.1250573049452193
$ws.Range("Q3").Value = 3.550711525519889
$ws.Range("R3").Value = 31.956403729679
$ws.Range("S3").Value = 0.03166695571928662
$ws.Range("T3").Value = 0.03166695571928661

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.9584923333333334
$ws.Range("H4").Value = 2.875477
$ws.Range("I4").Value = 0.2532195598902293
$ws.Range("J4").Value = 0.2532195598902293
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.371866333333333
$ws.Range("N4").Value = 19.115599
$ws.Range("O4").Value = 0.2151042422246107
$ws.Range("P4").Value = 0.2151042422246107
$ws.Range("Q4").Value = 6.107385029524778
$ws.Range("R4").Value = 54.966465265723
$ws.Range("S4").Value = 0.05446860154663721
$ws.Range("T4").Value = 0.0544686015466372

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.9584923333333334
$ws.Range("H5").Value = 2.875477
$ws.Range("I5").Value = 0.2532195598902293
$ws.Range("J5").Value = 0.2532195598902293
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 13.51475866666667
$ws.Range("N5").Value = 40.544276
$ws.Range("O5").Value = 0.4562371163741963
$ws.Range("P5").Value = 0.4562371163741963
$ws.Range("Q5").Value = 12.95379256885022
$ws.Range("R5").Value = 116.584133119652
$ws.Range("S5").Value = 0.1155281618138613
$ws.Range("T5").Value = 0.1155281618138613

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.9584923333333334
$ws.Range("H6").Value = 2.875477
$ws.Range("I6").Value = 0.2532195598902293
$ws.Range("J6").Value = 0.2532195598902293
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.6319533333333333
$ws.Range("N6").Value = 1.89586
$ws.Range("O6").Value = 0.02133375619900535
$ws.Range("P6").Value = 0.02133375619900535
$ws.Range("Q6").Value = 0.6057224250244444
$ws.Range("R6").Value = 5.451501825219999
$ws.Range("S6").Value = 0.005402124355517584
$ws.Range("T6").Value = 0.005402124355517584

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.82673
$ws.Range("H7").Value = 8.48019
$ws.Range("I7").Value = 0.7467804401097707
$ws.Range("J7").Value = 0.7467804401097707
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.399171333333334
$ws.Range("N7").Value = 16.197514
$ws.Range("O7").Value = 0.1822675802569684
$ws.Range("P7").Value = 0.1822675802569683
$ws.Range("Q7").Value = 15.26199958307333
$ws.Range("R7").Value = 137.35799624766
$ws.Range("S7").Value = 0.1361138638020418
$ws.Range("T7").Value = 0.1361138638020418

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.82673
$ws.Range("H8").Value = 8.48019
$ws.Range("I8").Value = 0.7467804401097707
$ws.Range("J8").Value = 0.7467804401097707
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.704475666666667
$ws.Range("N8").Value = 11.113427
$ws.Range("O8").Value = 0.1250573049452193
$ws.Range("P8").Value = 0.1250573049452193
$ws.Range("Q8").Value = 10.47155250123667
$ws.Range("R8").Value = 94.24397251113
$ws.Range("S8").Value = 0.09339034922593267
$ws.Range("T8").Value = 0.09339034922593266

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.82673
$ws.Range("H9").Value = 8.48019
$ws.Range("I9").Value = 0.7467804401097707
$ws.Range("J9").Value = 0.7467804401097707
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 6.371866333333333
$ws.Range("N9").Value = 19.115599
$ws.Range("O9").Value = 0.2151042422246107
$ws.Range("P9").Value = 0.2151042422246107
$ws.Range("Q9").Value = 18.01154572042333
$ws.Range("R9").Value = 162.10391148381
$ws.Range("S9").Value = 0.1606356406779735
$ws.Range("T9").Value = 0.1606356406779735

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.82673
$ws.Range("H10").Value = 8.48019
$ws.Range("I10").Value = 0.7467804401097707
$ws.Range("J10").Value = 0.7467804401097707
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 13.51475866666667
$ws.Range("N10").Value = 40.544276
$ws.Range("O10").Value = 0.4562371163741963
$ws.Range("P10").Value = 0.4562371163741963
$ws.Range("Q10").Value = 38.20257376582666
$ws.Range("R10").Value = 343.82316389244
$ws.Range("S10").Value = 0.340708954560335
$ws.Range("T10").Value = 0.340708954560335

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.82673
$ws.Range("H11").Value = 8.48019
$ws.Range("I11").Value = 0.7467804401097707
$ws.Range("J11").Value = 0.7467804401097707
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.6319533333333333
$ws.Range("N11").Value = 1.89586
$ws.Range("O11").Value = 0.02133375619900535
$ws.Range("P11").Value = 0.02133375619900535
$ws.Range("Q11").Value = 1.786361445933333
$ws.Range("R11").Value = 16.0772530134
$ws.Range("S11").Value = 0.01593163184348776
$ws.Range("T11").Value = 0.01593163184348776
